# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 18
$ws1.Range("F3").Value  = 1436
$ws1.Range("F7").Value  = 12003
$ws1.Range("F8").Value  = 4457
$ws1.Range("F9").Value  = 37
$ws1.Range("F10").Value = 56
$ws1.Range("F12").Value = 23
$ws1.Range("F13").Value = 2577
$ws1.Range("F17").Value = 5199
$ws1.Range("F20").Value = 541
$ws1.Range("F21").Value = 11403
$ws1.Range("F22").Value = 11425

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 18
$ws4.Range("F3").Value  = 1436
$ws4.Range("F7").Value  = 12003
$ws4.Range("F8").Value  = 4457
$ws4.Range("F9").Value  = 37
$ws4.Range("F10").Value = 56
$ws4.Range("F12").Value = 23
$ws4.Range("F13").Value = 2577
$ws4.Range("F18").Value = 5199
$ws4.Range("F21").Value = 541
$ws4.Range("F22").Value = 11403
$ws4.Range("F23").Value = 11425
